# Update column A course-id values on the active sheet so that a zero is
# inserted before the last digit of each id (making the ids "more readable"
# as integer-looking course codes rather than short numeric codes).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 21000
$ws.Range("A3").Value = 21100
$ws.Range("A4").Value = 21200
$ws.Range("A5").Value = 31300
$ws.Range("A6").Value = 31400
$ws.Range("A7").Value = 31500
$ws.Range("A8").Value = 33300
$ws.Range("A9").Value = 32200
$ws.Range("A10").Value = 44300
$ws.Range("A11").Value = 45100
$ws.Range("A12").Value = 45300
$ws.Range("A13").Value = 410001
$ws.Range("A14").Value = 410002
$ws.Range("A15").Value = 410003
$ws.Range("A16").Value = 33000
$ws.Range("A17").Value = 42200
$ws.Range("A18").Value = 44100
$ws.Range("A19").Value = 43600
$ws.Range("A20").Value = 43200
$ws.Range("A21").Value = 43100
$ws.Range("A22").Value = 429000
$ws.Range("A23").Value = 44500
$ws.Range("A24").Value = 41500
$ws.Range("A25").Value = 47100
$ws.Range("A26").Value = 47200
$ws.Range("A27").Value = 42500
$ws.Range("A28").Value = 41300
$ws.Range("A29").Value = 42000
$ws.Range("A30").Value = 47300
$ws.Range("A31").Value = 42300
$ws.Range("A32").Value = 46100
$ws.Range("A33").Value = 43300
$ws.Range("A34").Value = 43400
$ws.Range("A35").Value = 23101
$ws.Range("A36").Value = 23201
$ws.Range("A37").Value = 25101
$ws.Range("A38").Value = 25201
